$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Suicidal"
$ws.Range("A12").Value = "Suicidal"

[void]$ws.Range("A12").Select()
